$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "30.176.17"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  -3.45%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.919.96"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  -3.07%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9999"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  -0.32%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "246.10"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -2.82%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.7181"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -8.59%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.9986"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  -0.44%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3237"
$ws.Range("D8").ClearFormats()

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "26.26"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +2.62%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.06835"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -1.29%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.7926"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -4.65%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07919"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  -2.23%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.917.30"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  -3.59%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.374"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -1.17%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "94.27"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -5.86%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "14.39"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +4.07%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "259.06"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  -4.86%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "30.183.43"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -3.45%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "5.817"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +1.94%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.000007883"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -0.48%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "2.166.97"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -3.76%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.9983"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -0.50%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.000"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -0.30%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.836"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -1.20%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.623"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +0.15%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "160.20"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -2.65%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.1326"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -9.92%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.72"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -5.10%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.222"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +2.23%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.356"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -0.15%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.542"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -1.45%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.400"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -3.21%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.178"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -3.18%  "

$ws.Range("E34").Value = "  -2.67%  "

$ws.Range("E35").Value = "  -1.58%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7384"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -1.64%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.728"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -2.42%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01936"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -3.12%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.805"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -3.65%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "79.84"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +2.36%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "6.467"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -1.99%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.4403"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -4.73%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.006"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -2.28%  "

$ws.Range("E44").Value = "  -0.41%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.8293"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -2.60%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "102.16"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -3.03%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.707"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -2.61%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "7.236"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -3.09%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "35.95"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -1.25%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.4086"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -4.10%  "

$ws.Range("E51").Value = "  +2.90%  "
